$d = $word.ActiveDocument

# Change 1: underline "Praesent ornare fermentum turpis" within the first paragraph.
$r1 = $d.Content
$null = $r1.Find.Execute("Praesent ornare fermentum turpis", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r1.Font.Underline = 1

# Change 2: underline "posuere nec odio" (keeping existing bold) near the end.
$r2 = $d.Content
$null = $r2.Find.Execute("posuere nec odio", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r2.Font.Underline = 1
